# species.xlsx revision:
#  - rename header "Species" -> "species"
#  - normalize species names (underscore -> space) and "Non_native" -> "Non-native"
#  - re-sort data rows by species code (column B) and swap the
#    "Juncus bufonius / jubu" record for a new "Juncus gerardii / juge" record
#  - remove the yellow highlight fill used on the (old) Juncus bufonius row,
#    now carried by the Juncus gerardii row plus its "found as" cell
#  - move the active selection to A14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    1  = @("species", "code", "found as ", "endemic", "life", "root", "clade")
    2  = @("Achillea millefolium", "achmi", "seed", "Native", "perennial", "rhizomatous", "forb")
    3  = @("Agropyron repens", "agre", "vegetation", "Non-native", "perennial", "rhizomatous", "graminoid")
    4  = @("Agrostis spp.", "agsp", "both", "Non-native", "perennial", "rhizomatous", "graminoid")
    5  = @("Atriplex patula", "atpa", "both", "Non-native", "annual", "taproot", "forb")
    6  = @("Carex lyngbyei", "caly", "both", "Native", "perennial", "rhizomatous", "graminoid")
    7  = @("Cirsium sp", "cirsp", "seed", "Non-native", "perennial", "taproot", "forb")
    8  = @("Cotula coronopifolia", "coco", "both", "Non-native", "perennial", "fibrous", "forb")
    9  = @("Daucus carota", "daca", "vegetation", "Non-native", "perennial", "taproot", "forb")
    10 = @("Deschampsia cespitosa", "dece", "both", "Native", "perennial", "fibrous", "graminoid")
    11 = @("Distichlis spicata", "disp", "vegetation", "Native", "perennial", "rhizomatous", "graminoid")
    12 = @("Eleocharis parvula", "elpar", "both", "Native", "perennial", "rhizomatous", "graminoid")
    13 = @("Epilobium ciliatum", "epci", "seed", "Native", "perennial", "rhizomatous", "forb")
    14 = @("Glaux maritima", "glma", "both", "Native", "perennial", "rhizomatous", "forb")
    15 = @("Isolepis cernua", "isce", "seed", "Native", "annual", "fibrous", "graminoid")
    16 = @("Juncus articulatus", "juar", "seed", "Native", "perennial", "rhizomatous", "graminoid")
    17 = @("Juncus balticus", "juba", "both", "Native", "perennial", "rhizomatous", "forb")
    18 = @("Juncus ensifolius", "juen", "seed", "Native", "perennial", "rhizomatous", "graminoid")
    19 = @("Juncus gerardii", "juge", "seed", "Non-native", "perennial", "rhizomatous", "graminoid")
    20 = @("Polygonum fowlerii", "pofo", "vegetation", "Native", "annual", "taproot", "forb")
    21 = @("Potentilla anserina pacifica", "popa", "vegetation", "Native", "perennial", "rhizomatous", "forb")
    22 = @("Salicornia depressa", "sade", "both", "Native", "annual", "fibrous", "forb")
    23 = @("Spergularia canadensis", "spca", "both", "Native", "annual", "rhizomatous", "forb")
    24 = @("Symphotrichium subspicatum", "sysu", "both", "Native", "perennial", "rhizomatous", "forb")
    25 = @("Triglochin maritima", "trma", "vegetation", "Native", "perennial", "rhizomatous", "forb")
    26 = @("Trifolium wormskioldii", "trwo", "vegetation", "Native", "perennial", "taproot", "forb")
}

$cols = @("A", "B", "C", "D", "E", "F", "G")

foreach ($r in 1..26) {
    $rowValues = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $rowValues[$i]
    }
}

# The previously highlighted row (Juncus bufonius) carried a yellow fill on
# its first two cells; that row is now "Juncus gerardii" (row 19) and the
# highlight is cleared (the third, "found as", cell is touched too so the
# whole record shares one plain, unfilled style).
$ws.Range("C19").Interior.ColorIndex = 6
$ws.Range("A19:C19").Interior.Pattern = -4142

# Update the saved selection in the sheet view.
$ws.Range("A14").Select() | Out-Null
